$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new Price values are plain decimal numbers (e.g. "598.33"). Excel would
# auto-convert such strings to numeric cell values, but the source data keeps them
# as text (matching the untouched rows in this column). Force a text number format
# on exactly those cells first so the assigned values stay text, then restore the
# default style afterwards so no visible formatting change remains.
$textRng = $ws.Range("D5,D6,D7,D9,D11,D13,D17,D18,D21,D22,D23,D25,D26,D27,D29,D30,D31,D32,D34,D35,D36,D38,D39,D45,D48,D50,D51")
foreach ($area in $textRng.Areas) {
    $area.NumberFormat = "@"
}

$ws.Range('D2').Value = '73.020.03'
$ws.Range('E2').Value = '  +3.13%  '
$ws.Range('D3').Value = '3.985.94'
$ws.Range('E3').Value = '  +1.22%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '598.33'
$ws.Range('E5').Value = '  +12.03%  '
$ws.Range('D6').Value = '160.90'
$ws.Range('E6').Value = '  +8.85%  '
$ws.Range('D7').Value = '0.682'
$ws.Range('E7').Value = '  -0.60%  '
$ws.Range('E8').Value = '  -0.14%  '
$ws.Range('D9').Value = '0.749'
$ws.Range('E9').Value = '  +1.66%  '
$ws.Range('E10').Value = '  +1.96%  '
$ws.Range('D11').Value = '54.28'
$ws.Range('E11').Value = '  -1.56%  '
$ws.Range('E12').Value = '  +1.19%  '
$ws.Range('D13').Value = '10.94'
$ws.Range('E13').Value = '  +3.54%  '
$ws.Range('D14').Value = '4.625.09'
$ws.Range('E14').Value = '  +1.23%  '
$ws.Range('D15').Value = '3.989.01'
$ws.Range('E15').Value = '  +1.29%  '
$ws.Range('E16').Value = '  +9.57%  '
$ws.Range('D17').Value = '14.04'
$ws.Range('E17').Value = '  +1.76%  '
$ws.Range('D18').Value = '20.33'
$ws.Range('E18').Value = '  -0.86%  '
$ws.Range('E19').Value = '  +0.31%  '
$ws.Range('D20').Value = '72.762.16'
$ws.Range('E20').Value = '  +2.86%  '
$ws.Range('D21').Value = '436.10'
$ws.Range('E21').Value = '  +3.15%  '
$ws.Range('D22').Value = '4.82'
$ws.Range('E22').Value = '  +14.70%  '
$ws.Range('D23').Value = '95.91'
$ws.Range('E23').Value = '  -0.96%  '
$ws.Range('E24').Value = '  -4.53%  '
$ws.Range('D25').Value = '14.25'
$ws.Range('E25').Value = '  -0.83%  '
$ws.Range('D26').Value = '4.33'
$ws.Range('E26').Value = '  +14.61%  '
$ws.Range('D27').Value = '11.27'
$ws.Range('E27').Value = '  -0.57%  '
$ws.Range('E28').Value = '  +1.49%  '
$ws.Range('D29').Value = '10.38'
$ws.Range('E29').Value = '  -2.08%  '
$ws.Range('D30').Value = '36.29'
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('D31').Value = '7.82'
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('D32').Value = '13.73'
$ws.Range('E32').Value = '  +3.24%  '
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('D34').Value = '48.16'
$ws.Range('E34').Value = '  -4.92%  '
$ws.Range('D35').Value = '666.89'
$ws.Range('E35').Value = '  -2.13%  '
$ws.Range('D36').Value = '70.91'
$ws.Range('E36').Value = '  +9.06%  '
$ws.Range('D37').Value = '0.0₃0902'
$ws.Range('E37').Value = '  +10.71%  '
$ws.Range('D38').Value = '0.437'
$ws.Range('E38').Value = '  -0.14%  '
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('E40').Value = '  +5.19%  '
$ws.Range('E41').Value = '  -2.56%  '
$ws.Range('E42').Value = '  -0.76%  '
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('E44').Value = '  +2.01%  '
$ws.Range('D45').Value = '10.58'
$ws.Range('E45').Value = '  +6.20%  '
$ws.Range('E46').Value = '  +0.62%  '
$ws.Range('E47').Value = '  +3.33%  '
$ws.Range('D48').Value = '2.62'
$ws.Range('E48').Value = '  -2.12%  '
$ws.Range('D49').Value = '2.884.65'
$ws.Range('E49').Value = '  +9.72%  '
$ws.Range('D50').Value = '3.04'
$ws.Range('E50').Value = '  +1.56%  '
$ws.Range('D51').Value = '3.39'
$ws.Range('E51').Value = '  +4.55%  '

foreach ($area in $textRng.Areas) {
    $area.Style = "Normal"
}

